$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.48"
$ws.Range("E2").Value = "'0.27%"
$ws.Range("D3").Value = "'41.05"
$ws.Range("E3").Value = "'0.09%"
$ws.Range("D4").Value = "'5.145"
$ws.Range("E4").Value = "'-1.50%"
$ws.Range("D5").Value = "'0.07606"
$ws.Range("E5").Value = "'-0.68%"
$ws.Range("E6").Value = "'0.20%"
$ws.Range("D7").Value = "'1.678"
$ws.Range("E7").Value = "'2.19%"
$ws.Range("D8").Value = "'0.9352"
$ws.Range("E8").Value = "'2.26%"
$ws.Range("E9").Value = "'-2.88%"
$ws.Range("D10").Value = "'0.1820"
$ws.Range("E10").Value = "'-0.19%"
$ws.Range("D11").Value = "'0.09047"
$ws.Range("E11").Value = "'-1.27%"
$ws.Range("D12").Value = "'0.04137"
$ws.Range("E12").Value = "'0.89%"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'0.18%"
$ws.Range("D14").Value = "'0.001285"
$ws.Range("E14").Value = "'1.86%"
$ws.Range("D15").Value = "'0.005887"
$ws.Range("E15").Value = "'0.19%"
$ws.Range("D16").Value = "'3.339"
$ws.Range("E16").Value = "'-0.21%"
$ws.Range("D18").Value = "'0.3356"
$ws.Range("E18").Value = "'0.60%"
$ws.Range("D19").Value = "'7.594"
$ws.Range("E19").Value = "'1.90%"
$ws.Range("D20").Value = "'0.1341"
$ws.Range("E20").Value = "'-3.73%"
$ws.Range("D22").Value = "'0.03954"
$ws.Range("E22").Value = "'-2.23%"
$ws.Range("D23").Value = "'0.001279"
$ws.Range("E23").Value = "'1.25%"
$ws.Range("D24").Value = "'0.003981"
$ws.Range("E24").Value = "'-6.95%"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("E25").Value = "'6.16%"
$ws.Range("D26").Value = "'0.0003041"
$ws.Range("D38").Value = "'0.02427"
$ws.Range("E38").Value = "'-2.19%"
$ws.Range("D39").Value = "'0.05167"
$ws.Range("E39").Value = "'-3.22%"
$ws.Range("D40").Value = "'0.007721"
$ws.Range("E40").Value = "'-1.71%"
$ws.Range("D41").Value = "'0.1301"
$ws.Range("E41").Value = "'-0.80%"
$ws.Range("D42").Value = "'0.007590"
$ws.Range("E42").Value = "'15.25%"
$ws.Range("D43").Value = "'0.003301"
$ws.Range("E43").Value = "'72.49%"
$ws.Range("D44").Value = "'0.007576"
$ws.Range("E44").Value = "'-1.09%"
$ws.Range("D45").Value = "'0.3328"
$ws.Range("E45").Value = "'-0.74%"
$ws.Range("E46").Value = "'-2.33%"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("D48").Value = "'0.2709"
$ws.Range("E48").Value = "'36.59%"
$ws.Range("D49").Value = "'0.004202"
$ws.Range("E49").Value = "'35.23%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.16%"
